$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 135.4
$ws.Range("I2").Value = 120.625
$ws.Range("J2").Value = 194.5
$ws.Range("K2").Value = 120.625
$ws.Range("L2").Value = 194.5
$ws.Range("M2").Value = -7.625
$ws.Range("N2").Value = -420.5
# Row 4
$ws.Range("H4").Value = 343.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 29
$ws.Range("H29").Value = 324.75
$ws.Range("I29").Value = 324.75
$ws.Range("K29").Value = 974.25
$ws.Range("M29").Value = -693.25
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 196.16667
$ws.Range("I4").Value = 196.16667
$ws.Range("K4").Value = 196.16667
$ws.Range("M4").Value = -80.16667000000001
# Row 5
$ws.Range("H5").Value = 144
$ws.Range("I5").Value = 166
$ws.Range("K5").Value = 166
$ws.Range("M5").Value = -54
# Row 101
$ws.Range("H101").Value = 35000.5
$ws.Range("J101").Value = 35000.5
$ws.Range("L101").Value = 35000.5
$ws.Range("N101").Value = -41490.5
# Row 110
$ws.Range("H110").Value = 570.2308
$ws.Range("I110").Value = 550
$ws.Range("K110").Value = 550
$ws.Range("M110").Value = 1495
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 144
$ws.Range("I4").Value = 166
$ws.Range("K4").Value = 166
$ws.Range("M4").Value = -51
# Row 20
$ws.Range("H20").Value = 5543.6665
$ws.Range("I20").Value = 5152.4
$ws.Range("K20").Value = 5152.4
$ws.Range("M20").Value = -4905.4
# Row 22
$ws.Range("H22").Value = 419.6
$ws.Range("I22").Value = 419.6
$ws.Range("K22").Value = 419.6
$ws.Range("M22").Value = -246.6
# Row 86
$ws.Range("H86").Value = 3348.6667
$ws.Range("I86").Value = 3312.5386
$ws.Range("J86").Value = 3583.5
$ws.Range("K86").Value = 3312.5386
$ws.Range("L86").Value = 3583.5
$ws.Range("M86").Value = -2189.5386
$ws.Range("N86").Value = -5829.5
# Row 89
$ws.Range("H89").Value = 3348.6667
$ws.Range("I89").Value = 3312.5386
$ws.Range("J89").Value = 3583.5
$ws.Range("K89").Value = 16562.693
$ws.Range("L89").Value = 17917.5
$ws.Range("M89").Value = -10946.693
$ws.Range("N89").Value = -29149.5
# Row 103
$ws.Range("H103").Value = 8000
$ws.Range("J103").Value = 8000
$ws.Range("L103").Value = 8000
$ws.Range("N103").Value = -10344

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 108.833336
$ws.Range("I7").Value = 108.833336
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 108.833336
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 4.166663999999997
# Row 22
$ws.Range("H22").Value = 412
$ws.Range("I22").Value = 413.7143
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 413.7143
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -63.71429999999998
$ws.Range("N22").Value = -1100
# Row 31
$ws.Range("H31").Value = 2678.8667
$ws.Range("I31").Value = 1663.25
$ws.Range("K31").Value = 1663.25
$ws.Range("M31").Value = -1368.25
# Row 34
$ws.Range("H34").Value = 2678.8667
$ws.Range("I34").Value = 1663.25
$ws.Range("K34").Value = 1663.25
$ws.Range("M34").Value = -1461.25
# Row 58
$ws.Range("H58").Value = 1952.75
$ws.Range("I58").Value = 944.8
$ws.Range("J58").Value = 3632.6667
$ws.Range("K58").Value = 944.8
$ws.Range("L58").Value = 3632.6667
$ws.Range("M58").Value = -741.8
$ws.Range("N58").Value = -4038.6667
# Row 99
$ws.Range("H99").Value = 5011.647
$ws.Range("I99").Value = 4265.7
$ws.Range("K99").Value = 4265.7
$ws.Range("M99").Value = -2767.7
# Row 107
$ws.Range("H107").Value = 917.375
$ws.Range("I107").Value = 965.8333
$ws.Range("J107").Value = 772
$ws.Range("K107").Value = 965.8333
$ws.Range("L107").Value = 772
$ws.Range("M107").Value = 954.1667
$ws.Range("N107").Value = -4612
# Row 121
$ws.Range("H121").Value = 60300
$ws.Range("J121").Value = 60300
$ws.Range("L121").Value = 60300
$ws.Range("N121").Value = -62920
# Row 126
$ws.Range("H126").Value = 5011.647
$ws.Range("I126").Value = 4265.7
$ws.Range("K126").Value = 12797.1
$ws.Range("M126").Value = -10327.1
# Row 134
$ws.Range("H134").Value = 999
$ws.Range("I134").Value = 999
$ws.Range("K134").Value = 2997
$ws.Range("M134").Value = -462
# Row 136
$ws.Range("H136").Value = 1952.75
$ws.Range("I136").Value = 944.8
$ws.Range("J136").Value = 3632.6667
$ws.Range("K136").Value = 2834.4
$ws.Range("L136").Value = 10898.0001
$ws.Range("M136").Value = -284.3999999999996
$ws.Range("N136").Value = -15998.0001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1067.909
$ws.Range("J34").Value = 1363.8572
$ws.Range("L34").Value = 4091.5716
$ws.Range("N34").Value = -4259.571599999999
# Row 38
$ws.Range("H38").Value = 1009.75
$ws.Range("I38").Value = 3780.5
$ws.Range("J38").Value = 86.166664
$ws.Range("K38").Value = 11341.5
$ws.Range("L38").Value = 258.499992
$ws.Range("M38").Value = -10994.5
$ws.Range("N38").Value = -952.499992
# Row 39
$ws.Range("H39").Value = 2749.75
$ws.Range("J39").Value = 3166.3333
$ws.Range("L39").Value = 9498.999899999999
$ws.Range("N39").Value = -10086.9999
# Row 40
$ws.Range("H40").Value = 12
$ws.Range("I40").Value = 7.4
$ws.Range("J40").Value = 35
$ws.Range("K40").Value = 29.6
$ws.Range("L40").Value = 140
$ws.Range("M40").Value = 39.4
$ws.Range("N40").Value = -278
# Row 55
$ws.Range("H55").Value = 2999.5
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 3999
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 11997
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -12351
# Row 68
$ws.Range("H68").Value = 726.5
$ws.Range("I68").Value = 726.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2179.5
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -1368.5
# Row 71
$ws.Range("H71").Value = 726.5
$ws.Range("I71").Value = 726.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6538.5
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -2482.5
# Row 92
$ws.Range("H92").Value = 100
$ws.Range("I92").Value = 100
$ws.Range("K92").Value = 300
$ws.Range("M92").Value = 948

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8000
$ws.Range("J70").Value = 8000
$ws.Range("L70").Value = 8000
$ws.Range("N70").Value = -8540
# Row 73
$ws.Range("H73").Value = 8000
$ws.Range("J73").Value = 8000
$ws.Range("L73").Value = 8000
$ws.Range("N73").Value = -9872
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
# Row 132
$ws.Range("H132").Value = 1536.1666
$ws.Range("I132").Value = 1453.4
$ws.Range("K132").Value = 4360.200000000001
$ws.Range("M132").Value = -1830.200000000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 97
$ws.Range("H97").Value = 12855
$ws.Range("J97").Value = 12855
$ws.Range("L97").Value = 12855
$ws.Range("N97").Value = -14837
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0
# Row 122
$ws.Range("H122").Value = 9994.5
$ws.Range("I122").Value = 9994.5
$ws.Range("K122").Value = 29983.5
$ws.Range("M122").Value = -27533.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 72500
$ws.Range("J27").Value = 72500
$ws.Range("L27").Value = 72500
$ws.Range("N27").Value = -72638
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
# Row 107
$ws.Range("H107").Value = 847.25
$ws.Range("I107").Value = 795.3333
$ws.Range("J107").Value = 1003
$ws.Range("K107").Value = 2385.9999
$ws.Range("L107").Value = 3009
$ws.Range("M107").Value = -465.9998999999998
$ws.Range("N107").Value = -6849

Write-Host "Applied all Rafflesia_Profits updates"